$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.893.72'
$ws.Range('E2').Value = '  +1.86%  '

$ws.Range('D3').Value = '1.769.79'
$ws.Range('E3').Value = '  +2.37%  '

$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').Value = "'328.10"
$ws.Range('E5').Value = '  +2.00%  '

$ws.Range('E6').Value = '  -0.30%  '

$ws.Range('D7').Value = "'0.4483"
$ws.Range('E7').Value = '  -0.74%  '

$ws.Range('D8').Value = "'0.3559"
$ws.Range('E8').Value = '  +1.29%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'42.13"
$ws.Range('E9').Value = '  +1.58%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07430"

$ws.Range('D11').Value = "'1.103"
$ws.Range('E11').Value = '  +2.82%  '

$ws.Range('D12').Value = "'1.001"
$ws.Range('E12').Value = '  -0.23%  '

$ws.Range('D13').Value = "'20.99"
$ws.Range('E13').Value = '  +3.33%  '

$ws.Range('D14').Value = "'6.033"
$ws.Range('E14').Value = '  +2.40%  '

$ws.Range('D15').Value = "'7.251"
$ws.Range('E15').Value = '  +3.10%  '

$ws.Range('D16').Value = '1.769.62'
$ws.Range('E16').Value = '  +1.65%  '

$ws.Range('D17').Value = "'93.30"
$ws.Range('E17').Value = '  +2.43%  '

$ws.Range('D18').Value = "'0.00001063"
$ws.Range('E18').Value = '  +1.37%  '

$ws.Range('D19').Value = "'0.06432"
$ws.Range('E19').Value = '  +1.50%  '

$ws.Range('E20').Value = '  -0.37%  '

$ws.Range('D21').Value = "'17.12"
$ws.Range('E21').Value = '  +3.54%  '

$ws.Range('D22').Value = "'5.786"
$ws.Range('E22').Value = '  +1.01%  '

$ws.Range('D23').Value = '27.940.75'
$ws.Range('E23').Value = '  +1.87%  '

$ws.Range('E24').Value = '  +2.52%  '

$ws.Range('D25').Value = "'2.105"
$ws.Range('E25').Value = '  +0.84%  '

$ws.Range('D26').Value = "'162.15"
$ws.Range('E26').Value = '  +0.21%  '

$ws.Range('D27').Value = "'20.42"
$ws.Range('E27').Value = '  +3.32%  '

$ws.Range('D28').Value = '1.974.24'
$ws.Range('E28').Value = '  +2.06%  '

$ws.Range('D29').Value = "'2.169"
$ws.Range('E29').Value = '  +6.10%  '

$ws.Range('D30').Value = "'124.84"
$ws.Range('E30').Value = '  +0.53%  '

$ws.Range('D31').Value = "'1.108"
$ws.Range('E31').Value = '  +6.46%  '

$ws.Range('D32').Value = "'0.09208"
$ws.Range('E32').Value = '  +1.44%  '

$ws.Range('D33').Value = "'5.652"
$ws.Range('E33').Value = '  +6.04%  '

$ws.Range('D34').Value = "'3.672"
$ws.Range('E34').Value = '  +0.60%  '

$ws.Range('D35').Value = "'11.88"
$ws.Range('E35').Value = '  +2.44%  '

$ws.Range('E36').Value = '  +1.49%  '

$ws.Range('D37').Value = "'0.06102"
$ws.Range('E37').Value = '  +2.70%  '

$ws.Range('D38').Value = "'0.2104"
$ws.Range('E38').Value = '  +3.14%  '

$ws.Range('D39').Value = "'0.6337"
$ws.Range('E39').Value = '  +2.23%  '

$ws.Range('D40').Value = "'4.976"
$ws.Range('E40').Value = '  +3.00%  '

$ws.Range('D41').Value = "'1.184"
$ws.Range('E41').Value = '  -0.03%  '

$ws.Range('D42').Value = "'1.393"
$ws.Range('E42').Value = '  +1.75%  '

$ws.Range('D43').Value = "'7.918"
$ws.Range('E43').Value = '  +3.05%  '

$ws.Range('D44').Value = "'13.32"
$ws.Range('E44').Value = '  +2.65%  '

$ws.Range('D45').Value = "'3.736"
$ws.Range('E45').Value = '  +1.30%  '

$ws.Range('D46').Value = "'0.5896"
$ws.Range('E46').Value = '  +2.07%  '

$ws.Range('D47').Value = "'122.77"
$ws.Range('E47').Value = '  +0.91%  '

$ws.Range('D48').Value = "'1.962"
$ws.Range('E48').Value = '  +2.56%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = "'1.139"
$ws.Range('E49').Value = '  +3.08%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.06902"
$ws.Range('E50').Value = '  +1.27%  '

$ws.Range('D51').Value = "'73.11"
$ws.Range('E51').Value = '  +3.38%  '
